{"js": "// Append the new \"Encapsulation\" Q&A block to the end of the document,\n// following the same \"blank line separated\" pattern used by every other\n// topic in this cheat-sheet (question / blank / answer / blank, with a\n// blank-paragraph gap before the next question).\nconst newParagraphs = [\n  \"\",\n  \"\",\n  \"\",\n  \"\",\n  \"What is Encapsulation?\",\n  \"\",\n  \"Encapsulation is a fundamental principle of Object-Oriented Programming (OOP) that restricts direct access to an object\\u2019s data while allowing controlled access through methods. It helps protect the integrity of the data and prevents unintended modifications. In C# Encapsulation is implemented using modifiers such as private and public. Encapsulation ensures that an object\\u2019s internal state is hidden from the outside and only exposed through controlled mechanisms. \",\n  \"\",\n  \"Why do we use Encapsulation?\",\n  \"\",\n  \"Encapsulation simplifies debugging and modification of the class without affecting other parts of the program. It also encourages modular code design, making components easier to reuse. By using Encapsulation, we can enforce constraints on data and ensure that objects remain in a valid state.   \",\n  \"\",\n  \"How do we use Encapsulation in our code?\",\n  \"\",\n  \"In code, we use Encapsulation by making fields of a class private, and providing getter and setter methods (properties) to access and modify those fields. For example, aantalOVPaaltjes is private and can't be accessed directly from outside the class. Another example is getAantalOVPaaltjes and setAantalOVPaaltjes they are public properties, providing controlled access to the private field. The getter retrieves the value, and the setter updates the value. This way, Encapsulation ensures that the internal data of the object is protected and modified only through valid operations, providing data integrity and beter code structure.\",\n];\n\nlet anchor = context.document.body.paragraphs.getLast();\nfor (const text of newParagraphs) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n}\nawait context.sync();\n", "ps1": "# Append the new \"Encapsulation\" Q&A block to the end of the document,\n# following the same \"blank line separated\" pattern used by every other\n# topic in this cheat-sheet (question / blank / answer / blank, with a\n# blank-paragraph gap before the next question).\n$d = $word.ActiveDocument\n\n$newParagraphs = @(\n  \"\",\n  \"\",\n  \"\",\n  \"\",\n  \"What is Encapsulation?\",\n  \"\",\n  \"Encapsulation is a fundamental principle of Object-Oriented Programming (OOP) that restricts direct access to an object\u2019s data while allowing controlled access through methods. It helps protect the integrity of the data and prevents unintended modifications. In C# Encapsulation is implemented using modifiers such as private and public. Encapsulation ensures that an object\u2019s internal state is hidden from the outside and only exposed through controlled mechanisms. \",\n  \"\",\n  \"Why do we use Encapsulation?\",\n  \"\",\n  \"Encapsulation simplifies debugging and modification of the class without affecting other parts of the program. It also encourages modular code design, making components easier to reuse. By using Encapsulation, we can enforce constraints on data and ensure that objects remain in a valid state.   \",\n  \"\",\n  \"How do we use Encapsulation in our code?\",\n  \"\",\n  \"In code, we use Encapsulation by making fields of a class private, and providing getter and setter methods (properties) to access and modify those fields. For example, aantalOVPaaltjes is private and can't be accessed directly from outside the class. Another example is getAantalOVPaaltjes and setAantalOVPaaltjes they are public properties, providing controlled access to the private field. The getter retrieves the value, and the setter updates the value. This way, Encapsulation ensures that the internal data of the object is protected and modified only through valid operations, providing data integrity and beter code structure.\"\n)\n\nforeach ($text in $newParagraphs) {\n  $d.Paragraphs.Last.Range.InsertParagraphAfter()\n  if ($text -ne \"\") {\n    $d.Paragraphs.Last.Range.Text = $text\n  }\n}\n"}
